# TST: Test unnamed columns with index_col for Excel (gh-18792)
# Adds a new "Sheet4" at the end of the workbook containing a small table
# that mimics a pandas DataFrame written with an unnamed index column:
#
#        col1  col2
#   i1    a     x
#   i2    b     y
#
# The header row (col1/col2) and the index column (i1/i2) are bold,
# centered, top-aligned and boxed with a thin border - matching the
# formatting pandas/openpyxl applies to index & header cells.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet so it lands at the
# end of the tab strip (and becomes the active sheet, like a freshly
# added/selected sheet in Excel).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Sheet4"

# --- values -----------------------------------------------------------
$ws.Range("B1").Value = "col1"
$ws.Range("C1").Value = "col2"

$ws.Range("A2").Value = "i1"
$ws.Range("B2").Value = "a"
$ws.Range("C2").Value = "x"

$ws.Range("A3").Value = "i2"
$ws.Range("B3").Value = "b"
$ws.Range("C3").Value = "y"

# --- formatting ---------------------------------------------------------
# Header cells (B1:C1) and index cells (A2:A3) get bold Calibri 11
# (the workbook theme's minor/body font, themed text color 1),
# centered/top-aligned, with a thin box border around each cell.
function Set-HeaderStyle($rng) {
    $rng.Font.Name = "Calibri"
    $rng.Font.ThemeFont = 1            # xlThemeFontMinor
    $rng.Font.Size = 11
    $rng.Font.Bold = $true
    $rng.Font.ThemeColor = 1           # xlThemeColorDark1 (theme="1")
    $rng.HorizontalAlignment = -4108   # xlCenter
    $rng.VerticalAlignment = -4160     # xlTop
    $rng.Borders.LineStyle = 1         # xlContinuous / thin
}

Set-HeaderStyle($ws.Range("B1:C1"))
Set-HeaderStyle($ws.Range("A2:A3"))

# --- selection / view ----------------------------------------------------
$ws.Range("A1:C3").Select()
